$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9.62847638162277
$ws.Range("D2").Value = 4.889477451973894
$ws.Range("E2").Value = 12.5372904526273
$ws.Range("F2").Value = 30.06386191256673
$ws.Range("G2").Value = 3.626191705832189
$ws.Range("I2").Value = 26.77851708482337
$ws.Range("L2").Value = 9.461361224737427
$ws.Range("N2").Value = 18.99769123883711
$ws.Range("O2").Value = 26.27224278197653
$ws.Range("C3").Value = 9.616662997744289
$ws.Range("D3").Value = 4.904063324356646
$ws.Range("E3").Value = 12.50587338210526
$ws.Range("F3").Value = 29.61645569846443
$ws.Range("G3").Value = 3.629809073164635
$ws.Range("I3").Value = 26.52703300356168
$ws.Range("L3").Value = 9.46215684784668
$ws.Range("N3").Value = 18.40031508502701
$ws.Range("O3").Value = 25.99078460847669
$ws.Range("C4").Value = 9.611393688489017
$ws.Range("D4").Value = 4.913540692466086
$ws.Range("E4").Value = 12.48947518498017
$ws.Range("F4").Value = 29.34776987288977
$ws.Range("G4").Value = 3.632147283375156
$ws.Range("I4").Value = 26.37914352774852
$ws.Range("L4").Value = 9.464457197895866
$ws.Range("N4").Value = 18.02485520896362
$ws.Range("O4").Value = 25.82415151807121
$ws.Range("C5").Value = 9.609746457278018
$ws.Range("D5").Value = 4.917534230268339
$ws.Range("E5").Value = 12.48352350784842
$ws.Range("F5").Value = 29.23993525346501
$ws.Range("G5").Value = 3.633129687416996
$ws.Range("I5").Value = 26.32057785934654
$ws.Range("L5").Value = 9.465850124014667
$ws.Range("N5").Value = 17.86990355188765
$ws.Range("O5").Value = 25.75787497046569
$ws.Range("C6").Value = 9.609503163802685
$ws.Range("D6").Value = 4.918205300668256
$ws.Range("E6").Value = 12.4825794579916
$ws.Range("F6").Value = 29.22213357267411
$ws.Range("G6").Value = 3.633294603725532
$ws.Range("I6").Value = 26.3109574692815
$ws.Range("L6").Value = 9.466108923981695
$ws.Range("N6").Value = 17.84406337566597
$ws.Range("O6").Value = 25.74697017307625
$ws.Range("C7").Value = 9.611369447426789
$ws.Range("D7").Value = 4.913594018198795
$ws.Range("E7").Value = 12.48939195584017
$ws.Range("F7").Value = 29.34630868273597
$ws.Range("G7").Value = 3.632160412569012
$ws.Range("I7").Value = 26.37834672866202
$ws.Range("L7").Value = 9.46447413928604
$ws.Range("N7").Value = 18.02277304767603
$ws.Range("O7").Value = 25.82325100505323
$ws.Range("C8").Value = 9.623992090321494
$ws.Range("D8").Value = 4.894398571535359
$ws.Range("E8").Value = 12.52585992806409
$ws.Range("F8").Value = 29.90843246479154
$ws.Range("G8").Value = 3.627414729997339
$ws.Range("I8").Value = 26.69048884695097
$ws.Range("L8").Value = 9.461259512296644
$ws.Range("N8").Value = 18.79364780656867
$ws.Range("O8").Value = 26.17395749425654
$ws.Range("C9").Value = 9.664425100048867
$ws.Range("D9").Value = 4.860882506918918
$ws.Range("E9").Value = 12.62013471018199
$ws.Range("F9").Value = 31.0521652329615
$ws.Range("G9").Value = 3.619032815269698
$ws.Range("I9").Value = 27.35153731079814
$ws.Range("L9").Value = 9.469328790192943
$ws.Range("N9").Value = 20.2273683202997
$ws.Range("O9").Value = 26.90738542557753
$ws.Range("C10").Value = 9.703584947228398
$ws.Range("D10").Value = 4.83875671749872
$ws.Range("E10").Value = 12.70298026706912
$ws.Range("F10").Value = 31.90909069558782
$ws.Range("G10").Value = 3.613431099136563
$ws.Range("I10").Value = 27.86307116331835
$ws.Range("L10").Value = 9.484005659677701
$ws.Range("N10").Value = 21.22223697909767
$ws.Range("O10").Value = 27.46939495718158
$ws.Range("C11").Value = 9.723422804829298
$ws.Range("D11").Value = 4.829230182980272
$ws.Range("E11").Value = 12.74354063641795
$ws.Range("F11").Value = 32.30072432240137
$ws.Range("G11").Value = 3.611002066040073
$ws.Range("I11").Value = 28.10050732186091
$ws.Range("L11").Value = 9.492575029731363
$ws.Range("N11").Value = 21.66018057919901
$ws.Range("O11").Value = 27.7290626869194
$ws.Range("C12").Value = 9.731222780033178
$ws.Range("D12").Value = 4.825699937287694
$ws.Range("E12").Value = 12.75930534164641
$ws.Range("F12").Value = 32.44913868872323
$ws.Range("G12").Value = 3.610099282081451
$ws.Range("I12").Value = 28.19102285860927
$ws.Range("L12").Value = 9.496091200808445
$ws.Range("N12").Value = 21.82377585682186
$ws.Range("O12").Value = 27.82788060060437
$ws.Range("C13").Value = 9.72953017036326
$ws.Range("D13").Value = 4.82645680645001
$ws.Range("E13").Value = 12.75589223281854
$ws.Range("F13").Value = 32.41717239182119
$ws.Range("G13").Value = 3.610292956648079
$ws.Range("I13").Value = 28.17150311226839
$ws.Range("L13").Value = 9.495321890515607
$ws.Range("N13").Value = 21.78864458690801
$ws.Range("O13").Value = 27.80657812986101
$ws.Range("C14").Value = 9.724058765028046
$ws.Range("D14").Value = 4.828938200548297
$ws.Range("E14").Value = 12.74482953108404
$ws.Range("F14").Value = 32.3129330752967
$ws.Range("G14").Value = 3.610927452589677
$ws.Range("I14").Value = 28.10794243819736
$ws.Range("L14").Value = 9.492858880012459
$ws.Range("N14").Value = 21.67368539489659
$ws.Range("O14").Value = 27.73718323621465
$ws.Range("C15").Value = 9.720744751070228
$ws.Range("D15").Value = 4.830468179931631
$ws.Range("E15").Value = 12.7381058548877
$ws.Range("F15").Value = 32.24909359555794
$ws.Range("G15").Value = 3.611318315440208
$ws.Range("I15").Value = 28.06908598802562
$ws.Range("L15").Value = 9.491385490258972
$ws.Range("N15").Value = 21.60297336126124
$ws.Range("O15").Value = 27.69473769911087
$ws.Range("C16").Value = 9.702328954494435
$ws.Range("D16").Value = 4.839390119283604
$ws.Range("E16").Value = 12.70038671921722
$ws.Range("F16").Value = 31.88352013745232
$ws.Range("G16").Value = 3.613592232106549
$ws.Range("I16").Value = 27.84764285588846
$ws.Range("L16").Value = 9.483483634650071
$ws.Range("N16").Value = 21.19330956972086
$ws.Range("O16").Value = 27.45249809503873
$ws.Range("C17").Value = 9.691547637645566
$ws.Range("D17").Value = 4.845001240094928
$ws.Range("E17").Value = 12.67797786151845
$ws.Range("F17").Value = 31.6596095785052
$ws.Range("G17").Value = 3.615017665516506
$ws.Range("I17").Value = 27.71295294107361
$ws.Range("L17").Value = 9.479120216899572
$ws.Range("N17").Value = 20.93814219015166
$ws.Range("O17").Value = 27.30485499789356
$ws.Range("C18").Value = 9.685537218517799
$ws.Range("D18").Value = 4.848279315373632
$ws.Range("E18").Value = 12.66535985921963
$ws.Range("F18").Value = 31.53100161992067
$ws.Range("G18").Value = 3.615848764499717
$ws.Range("I18").Value = 27.63593341792934
$ws.Range("L18").Value = 9.476788728982751
$ws.Range("N18").Value = 20.79000725568362
$ws.Range("O18").Value = 27.22031705978937
$ws.Range("C19").Value = 9.683535040637615
$ws.Range("D19").Value = 4.84939793065706
$ws.Range("E19").Value = 12.66113438942691
$ws.Range("F19").Value = 31.48749262216938
$ws.Range("G19").Value = 3.616132091852769
$ws.Range("I19").Value = 27.60993566520034
$ws.Range("L19").Value = 9.476029967427641
$ws.Range("N19").Value = 20.73962067985785
$ws.Range("O19").Value = 27.19176236151834
$ws.Range("C20").Value = 9.692675611766626
$ws.Range("D20").Value = 4.844398680275358
$ws.Range("E20").Value = 12.68033533196438
$ws.Range("F20").Value = 31.68342776180243
$ws.Range("G20").Value = 3.614864764470275
$ws.Range("I20").Value = 27.72724484181598
$ws.Range("L20").Value = 9.479566270463517
$ws.Range("N20").Value = 20.96544799484619
$ws.Range("O20").Value = 27.32053291759324
$ws.Range("C21").Value = 9.725658065717594
$ws.Range("D21").Value = 4.828207260086014
$ws.Range("E21").Value = 12.748067978664
$ws.Range("F21").Value = 32.34354884325683
$ws.Range("G21").Value = 3.610740624174252
$ws.Range("I21").Value = 28.12659598749733
$ws.Range("L21").Value = 9.493574976728615
$ws.Range("N21").Value = 21.70751365554066
$ws.Range("O21").Value = 27.75755367790951
$ws.Range("C22").Value = 9.748889841475638
$ws.Range("D22").Value = 4.818075393189917
$ws.Range("E22").Value = 12.79469392360499
$ws.Range("F22").Value = 32.77555802606521
$ws.Range("G22").Value = 3.608144523164056
$ws.Range("I22").Value = 28.39107715453967
$ws.Range("L22").Value = 9.504310230473457
$ws.Range("N22").Value = 22.17935961385674
$ws.Range("O22").Value = 28.04597191127845
$ws.Range("C23").Value = 9.73633841889208
$ws.Range("D23").Value = 4.823441836119838
$ws.Range("E23").Value = 12.76959570744089
$ws.Range("F23").Value = 32.54498168427385
$ws.Range("G23").Value = 3.609521063088993
$ws.Range("I23").Value = 28.24962569373607
$ws.Range("L23").Value = 9.498436476824001
$ws.Range("N23").Value = 21.92877110911181
$ws.Range("O23").Value = 27.89181073488903
$ws.Range("C24").Value = 9.692165069142856
$ws.Range("D24").Value = 4.844670935041291
$ws.Range("E24").Value = 12.67926869310599
$ws.Range("F24").Value = 31.67265917207502
$ws.Range("G24").Value = 3.61493385489731
$ws.Range("I24").Value = 27.72078217329675
$ws.Range("L24").Value = 9.479364057936088
$ws.Range("N24").Value = 20.95310750188673
$ws.Range("O24").Value = 27.31344385160817
$ws.Range("C25").Value = 9.651817675125148
$ws.Range("D25").Value = 4.86950954277091
$ws.Range("E25").Value = 12.59222161091895
$ws.Range("F25").Value = 30.73921294634416
$ws.Range("G25").Value = 3.621202115862118
$ws.Range("I25").Value = 27.16786770945205
$ws.Range("L25").Value = 9.465607692995313
$ws.Range("N25").Value = 19.84905939529497
$ws.Range("O25").Value = 26.70455883351522
